$wb = $excel.ActiveWorkbook

# --- Sheet2: autofit column M (new WC_mean_FST mirror column) and update view ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Columns.Item(13).AutoFit()
$ws2.Range("F29").Select()

# --- Sheet3 ("Sample" sheet): swap the two FST tables, add titles above each ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.UsedRange.ClearContents()

$ws3.Range("A1").Value = "WC FST "
$ws3.Range("A2").Value = "Pop"
$ws3.Range("B2").Value = "SU18"
$ws3.Range("C2").Value = "AR18"
$ws3.Range("D2").Value = "SI18"
$ws3.Range("E2").Value = "SD18"
$ws3.Range("F2").Value = "YS21"
$ws3.Range("G2").Value = "IS21"
$ws3.Range("H2").Value = "UH21"
$ws3.Range("A3").Value = "SU18"
$ws3.Range("B3").Value = "NA"
$ws3.Range("C3").Value = "NA"
$ws3.Range("D3").Value = "NA"
$ws3.Range("E3").Value = "NA"
$ws3.Range("F3").Value = "NA"
$ws3.Range("G3").Value = "NA"
$ws3.Range("H3").Value = "NA"
$ws3.Range("A4").Value = "AR18"
$ws3.Range("B4").Value = 0.0027492240000000002
$ws3.Range("C4").Value = "NA"
$ws3.Range("D4").Value = "NA"
$ws3.Range("E4").Value = "NA"
$ws3.Range("F4").Value = "NA"
$ws3.Range("G4").Value = "NA"
$ws3.Range("H4").Value = "NA"
$ws3.Range("A5").Value = "SI18"
$ws3.Range("B5").Value = 0.0078321199999999997
$ws3.Range("C5").Value = 0.0047663849999999997
$ws3.Range("D5").Value = "NA"
$ws3.Range("E5").Value = "NA"
$ws3.Range("F5").Value = "NA"
$ws3.Range("G5").Value = "NA"
$ws3.Range("H5").Value = "NA"
$ws3.Range("A6").Value = "SD18"
$ws3.Range("B6").Value = 0.0033005629999999998
$ws3.Range("C6").Value = 0.0040624279999999999
$ws3.Range("D6").Value = 0.0059278723000000004
$ws3.Range("E6").Value = "NA"
$ws3.Range("F6").Value = "NA"
$ws3.Range("G6").Value = "NA"
$ws3.Range("H6").Value = "NA"
$ws3.Range("A7").Value = "YS21"
$ws3.Range("B7").Value = 0.005832439
$ws3.Range("C7").Value = 0.0042699840000000001
$ws3.Range("D7").Value = 0.00029826160000000001
$ws3.Range("E7").Value = 0.0032943080000000001
$ws3.Range("F7").Value = "NA"
$ws3.Range("G7").Value = "NA"
$ws3.Range("H7").Value = "NA"
$ws3.Range("A8").Value = "IS21"
$ws3.Range("B8").Value = 0.006783323
$ws3.Range("C8").Value = 0.005118175
$ws3.Range("D8").Value = 0.0045439466999999999
$ws3.Range("E8").Value = 0.0043807760000000003
$ws3.Range("F8").Value = 0.0009208749
$ws3.Range("G8").Value = "NA"
$ws3.Range("H8").Value = "NA"
$ws3.Range("A9").Value = "UH21"
$ws3.Range("B9").Value = 0.045215666000000002
$ws3.Range("C9").Value = 0.032584531999999999
$ws3.Range("D9").Value = 0.021307997499999998
$ws3.Range("E9").Value = 0.041324086000000003
$ws3.Range("F9").Value = 0.0133642404
$ws3.Range("G9").Value = 0.030229550000000001
$ws3.Range("H9").Value = "NA"
$ws3.Range("A11").Value = "p-value with 1000 bootstrap"
$ws3.Range("A12").Value = "POP"
$ws3.Range("B12").Value = "SU18"
$ws3.Range("C12").Value = "AR18"
$ws3.Range("D12").Value = "SI18"
$ws3.Range("E12").Value = "SD18"
$ws3.Range("F12").Value = "YS21"
$ws3.Range("G12").Value = "IS21"
$ws3.Range("H12").Value = "UH21"
$ws3.Range("A13").Value = "SU18"
$ws3.Range("B13").Value = "NA"
$ws3.Range("C13").Value = "NA"
$ws3.Range("D13").Value = "NA"
$ws3.Range("E13").Value = "NA"
$ws3.Range("F13").Value = "NA"
$ws3.Range("G13").Value = "NA"
$ws3.Range("H13").Value = "NA"
$ws3.Range("A14").Value = "AR18"
$ws3.Range("B14").Value = 0.13400000000000001
$ws3.Range("C14").Value = "NA"
$ws3.Range("D14").Value = "NA"
$ws3.Range("E14").Value = "NA"
$ws3.Range("F14").Value = "NA"
$ws3.Range("G14").Value = "NA"
$ws3.Range("H14").Value = "NA"
$ws3.Range("A15").Value = "SI18"
$ws3.Range("B15").Value = 0.0070000000000000001
$ws3.Range("C15").Value = 0.049000000000000002
$ws3.Range("D15").Value = "NA"
$ws3.Range("E15").Value = "NA"
$ws3.Range("F15").Value = "NA"
$ws3.Range("G15").Value = "NA"
$ws3.Range("H15").Value = "NA"
$ws3.Range("A16").Value = "SD18"
$ws3.Range("B16").Value = 0.106
$ws3.Range("C16").Value = 0.10199999999999999
$ws3.Range("D16").Value = 0.027
$ws3.Range("E16").Value = "NA"
$ws3.Range("F16").Value = "NA"
$ws3.Range("G16").Value = "NA"
$ws3.Range("H16").Value = "NA"
$ws3.Range("A17").Value = "YS21"
$ws3.Range("B17").Value = 0.010999999999999999
$ws3.Range("C17").Value = 0.029000000000000001
$ws3.Range("D17").Value = 0.46800000000000003
$ws3.Range("E17").Value = 0.10199999999999999
$ws3.Range("F17").Value = "NA"
$ws3.Range("G17").Value = "NA"
$ws3.Range("H17").Value = "NA"
$ws3.Range("A18").Value = "IS21"
$ws3.Range("B18").Value = 0.0040000000000000001
$ws3.Range("C18").Value = 0.014999999999999999
$ws3.Range("D18").Value = 0.042000000000000003
$ws3.Range("E18").Value = 0.083000000000000004
$ws3.Range("F18").Value = 0.374
$ws3.Range("G18").Value = "NA"
$ws3.Range("H18").Value = "NA"
$ws3.Range("A19").Value = "UH21"
$ws3.Range("B19").Value = 0
$ws3.Range("C19").Value = 0
$ws3.Range("D19").Value = 0
$ws3.Range("E19").Value = 0
$ws3.Range("F19").Value = 0
$ws3.Range("G19").Value = 0
$ws3.Range("H19").Value = "NA"

$ws3.Range("R12").Select()
